$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Import" row: trim the POJO import list down to just Insured.
$ws.Range("B2").Value = "com.redhat.prudential_poc.pojo.Insured"

# ACTION column for the four age-bracket rules: "免體檢" -> "PASS".
$ws.Range("E10").Value = '"PASS"'
$ws.Range("E11").Value = '"PASS"'
$ws.Range("E12").Value = '"PASS"'
$ws.Range("E13").Value = '"PASS"'

# Rows 12-13 (columns B:D) were missing the formatting already used by rows
# 10-11 in that block; copy it across so all four rule rows match.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B12:D13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Move the saved selection/active cell to B4.
$ws.Range("B4").Select() | Out-Null
